$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" stat (was Strike#, now K). Regenerated values per
# the updated save_data calc/write of s_vals for rows 2-11.
$newK = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    9  = 3
    10 = 4
    11 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
